# Apply updated "想去人数" (want-to-go count) figures to the 苏州-漫展信息 workbook.
# Three worksheets are affected: 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4).
# 本地生活 (sheet3) only has a header row and is not touched.

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12989
$ws1.Range("F4").Value = 29
$ws1.Range("F5").Value = 83
$ws1.Range("F6").Value = 93
$ws1.Range("F8").Value = 25
$ws1.Range("F10").Value = 12972
$ws1.Range("F12").Value = 43
$ws1.Range("F13").Value = 8716
$ws1.Range("F14").Value = 7727
$ws1.Range("F22").Value = 382
$ws1.Range("F24").Value = 331

# --- 演出 -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 2

# --- 全部类型 -------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12989
$ws4.Range("F5").Value = 29
$ws4.Range("F6").Value = 83
$ws4.Range("F7").Value = 93
$ws4.Range("F9").Value = 25
$ws4.Range("F11").Value = 12972
$ws4.Range("F13").Value = 43
$ws4.Range("F14").Value = 8716
$ws4.Range("F15").Value = 7727
$ws4.Range("F23").Value = 2
$ws4.Range("F25").Value = 382
$ws4.Range("F27").Value = 331
